# Update countries & provincias Spain
# Refresh the "Pais" sheet with the latest scrape snapshot (13:00 -> 14:17):
#  - bump the "datos actualizados" timestamp in A1
#  - a handful of neighbouring countries swapped rank order since the
#    previous scrape (Kazajistan/Nepal, Tunez/Dinamarca, Noruega/Zambia,
#    Seychelles/Brunei/Liechtenstein) so those rows get both the new
#    country label and the figures that now belong to that rank
#  - refresh the numeric columns (Casos totales, Nuevos casos, Casos
#    activos, Recuperados, Casos criticos, Muertes hoy, Muertes) for the
#    rows whose counts moved between the two snapshots

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("Pais")

# Timestamp banner
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 12 de Octubre de 2020 a las 14:17"

# Alemania (row 25)
$ws.Cells.Item(25, 2).Value = 327725
$ws.Cells.Item(25, 3).Value = 1434
$ws.Cells.Item(25, 5).Value = 41121
$ws.Cells.Item(25, 7).Value = 2
$ws.Cells.Item(25, 8).Value = 9704

# Kuwait (row 41)
$ws.Cells.Item(41, 2).Value = 111893
$ws.Cells.Item(41, 3).Value = 777
$ws.Cells.Item(41, 4).Value = 103802
$ws.Cells.Item(41, 5).Value = 7427
$ws.Cells.Item(41, 7).Value = 6
$ws.Cells.Item(41, 8).Value = 664

# Row 42: was Kazajistan, now Nepal (rank swap)
$ws.Cells.Item(42, 1).Value = "Nepal"
$ws.Cells.Item(42, 2).Value = 111802
$ws.Cells.Item(42, 3).Value = 4047
$ws.Cells.Item(42, 4).Value = 77277
$ws.Cells.Item(42, 5).Value = 33880
$ws.Cells.Item(42, 7).Value = 9
$ws.Cells.Item(42, 8).Value = 645

# Row 43: was Nepal, now Kazajistan (rank swap)
$ws.Cells.Item(43, 1).Value = "Kazajistan"
$ws.Cells.Item(43, 2).Value = 108831
$ws.Cells.Item(43, 3).Value = 74
$ws.Cells.Item(43, 4).Value = 104041
$ws.Cells.Item(43, 5).Value = 3044
$ws.Cells.Item(43, 8).Value = 1746

# Bielorrusia (row 55)
$ws.Cells.Item(55, 2).Value = 83998
$ws.Cells.Item(55, 3).Value = 464
$ws.Cells.Item(55, 4).Value = 77423
$ws.Cells.Item(55, 5).Value = 5674
$ws.Cells.Item(55, 7).Value = 5
$ws.Cells.Item(55, 8).Value = 901

# Suiza (row 58)
$ws.Cells.Item(58, 5).Value = 13944
$ws.Cells.Item(58, 7).Value = 4
$ws.Cells.Item(58, 8).Value = 2092

# Row 78: was Tunez, now Dinamarca (rank swap)
$ws.Cells.Item(78, 1).Value = "Dinamarca"
$ws.Cells.Item(78, 2).Value = 32811
$ws.Cells.Item(78, 3).Value = 389
$ws.Cells.Item(78, 4).Value = 26741
$ws.Cells.Item(78, 5).Value = 5399
$ws.Cells.Item(78, 7).Value = 2
$ws.Cells.Item(78, 8).Value = 671

# Row 79: was Dinamarca, now Tunez (rank swap)
$ws.Cells.Item(79, 1).Value = "Tunez"
$ws.Cells.Item(79, 2).Value = 32556
$ws.Cells.Item(79, 4).Value = 5032
$ws.Cells.Item(79, 5).Value = 27046
$ws.Cells.Item(79, 8).Value = 478

# Bosnia y Herzegovina (row 80)
$ws.Cells.Item(80, 2).Value = 30837
$ws.Cells.Item(80, 3).Value = 190
$ws.Cells.Item(80, 4).Value = 23587
$ws.Cells.Item(80, 5).Value = 6314
$ws.Cells.Item(80, 7).Value = 8
$ws.Cells.Item(80, 8).Value = 936

# Madagascar (row 93)
$ws.Cells.Item(93, 2).Value = 16726
$ws.Cells.Item(93, 3).Value = 8
$ws.Cells.Item(93, 4).Value = 16081
$ws.Cells.Item(93, 5).Value = 408

# Row 95: was Noruega, now Zambia (rank swap)
$ws.Cells.Item(95, 1).Value = "Zambia"
$ws.Cells.Item(95, 2).Value = 15549
$ws.Cells.Item(95, 3).Value = 91
$ws.Cells.Item(95, 4).Value = 14682
$ws.Cells.Item(95, 5).Value = 522
$ws.Cells.Item(95, 7).Value = 8
$ws.Cells.Item(95, 8).Value = 345

# Row 96: was Zambia, now Noruega (rank swap)
$ws.Cells.Item(96, 1).Value = "Noruega"
$ws.Cells.Item(96, 2).Value = 15524
$ws.Cells.Item(96, 4).Value = 11863
$ws.Cells.Item(96, 5).Value = 3385
$ws.Cells.Item(96, 7).Value = 1
$ws.Cells.Item(96, 8).Value = 276

# Islandia (row 145)
$ws.Cells.Item(145, 2).Value = 3582
$ws.Cells.Item(145, 3).Value = 56
$ws.Cells.Item(145, 4).Value = 2550
$ws.Cells.Item(145, 5).Value = 1022

# Vietnam (row 168)
$ws.Cells.Item(168, 2).Value = 1110
$ws.Cells.Item(168, 3).Value = 1
$ws.Cells.Item(168, 4).Value = 1025

# Row 194: was Seychelles, now Liechtenstein (rank rotation)
$ws.Cells.Item(194, 1).Value = "Liechtenstein"
$ws.Cells.Item(194, 4).Value = 124
$ws.Cells.Item(194, 5).Value = 23
$ws.Cells.Item(194, 8).Value = 1

# Row 195: was Brunei, now Seychelles (rank rotation)
$ws.Cells.Item(195, 1).Value = "Seychelles"
$ws.Cells.Item(195, 2).Value = 148
$ws.Cells.Item(195, 4).Value = 144
$ws.Cells.Item(195, 5).Value = 4
$ws.Cells.Item(195, 8).Value = 0

# Row 196: was Liechtenstein, now Brunei (rank rotation)
$ws.Cells.Item(196, 1).Value = "Brunei"
$ws.Cells.Item(196, 2).Value = 146
$ws.Cells.Item(196, 4).Value = 143
$ws.Cells.Item(196, 5).Value = 0
$ws.Cells.Item(196, 8).Value = 3
